$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123; this shifts existing rows 123:169 down to 124:170
$ws.Rows("123:123").Insert()

# Populate the newly inserted row 123 with the new data record
$ws.Cells.Item(123, 1).Value = 6
$ws.Cells.Item(123, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(123, 3).Value = 'Metropolitana'
$ws.Cells.Item(123, 4).Value = 44609
$ws.Cells.Item(123, 5).Value = 13
$ws.Cells.Item(123, 6).Value = 'Fruta'
$ws.Cells.Item(123, 7).Value = 100101
$ws.Cells.Item(123, 8).Value = 'Berries'
$ws.Cells.Item(123, 9).Value = 100101004
$ws.Cells.Item(123, 10).Value = 'Frambuesa'
$ws.Cells.Item(123, 11).Value = 'Sin especificar'
$ws.Cells.Item(123, 12).Value = 'Especial'
$ws.Cells.Item(123, 13).Value = 300
$ws.Cells.Item(123, 14).Value = 8000
$ws.Cells.Item(123, 15).Value = 8000
$ws.Cells.Item(123, 16).Value = 8000
$ws.Cells.Item(123, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(123, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(123, 19).Value = 4000
$ws.Cells.Item(123, 20).Value = 2
